$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "(""Valkyrie's Sword"", ['{1}{W}', 'Artifact — Equipment', 'When Valkyrie’s Sword enters the battlefield, you may pay {4}{W}. If you do, create a 4/4 white Angel Warrior creature token with flying and vigilance, then attach Valkyrie’s Sword to it.', 'Equipped creature gets +2/+1.', 'Equip {3}'])"
